# Turns the single-table "Test Document with Table" into the
# "Complex Test Document" with two tables, per the commit diff.

$d = $word.ActiveDocument

# --- 1. Text edits -------------------------------------------------------
$d.Content.Find.Execute("Test Document with Table", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Complex Test Document", 2)

$d.Content.Find.Execute("This is a test document.", $true, $false, $false, `
    $false, $false, $true, 1, $false, "This document has multiple tables.", 2)

# --- 2. Drop the old (single) table; we'll recreate it, plus a brand new
#        "Key/Value" table above it and a connecting paragraph, with exact
#        XML control over tblPr/tblLook so the markup matches Word's own
#        output instead of the engine's generic Tables.Add defaults.
$oldTable = $d.Tables(1)
$oldTable.Delete()

$introPara = $d.Paragraphs(2)
$insertAt = $introPara.Range.End - 1
$rng = $d.Range($insertAt, $insertAt)

$W = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = @"
<w:tbl $W>
<w:tblPr>
<w:tblW w:type='auto' w:w='0'/>
<w:tblLook w:firstColumn='1' w:firstRow='1' w:lastColumn='0' w:lastRow='0' w:noHBand='0' w:noVBand='1' w:val='04A0'/>
</w:tblPr>
<w:tblGrid>
<w:gridCol w:w='4320'/>
<w:gridCol w:w='4320'/>
</w:tblGrid>
<w:tr>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='4320'/></w:tcPr><w:p><w:r><w:t>Key</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='4320'/></w:tcPr><w:p><w:r><w:t>Value</w:t></w:r></w:p></w:tc>
</w:tr>
<w:tr>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='4320'/></w:tcPr><w:p><w:r><w:t>Status</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='4320'/></w:tcPr><w:p><w:r><w:t>Active</w:t></w:r></w:p></w:tc>
</w:tr>
</w:tbl>
<w:p $W><w:r><w:t>Here is another table:</w:t></w:r></w:p>
<w:tbl $W>
<w:tblPr>
<w:tblW w:type='auto' w:w='0'/>
<w:tblLook w:firstColumn='1' w:firstRow='1' w:lastColumn='0' w:lastRow='0' w:noHBand='0' w:noVBand='1' w:val='04A0'/>
</w:tblPr>
<w:tblGrid>
<w:gridCol w:w='2880'/>
<w:gridCol w:w='2880'/>
<w:gridCol w:w='2880'/>
</w:tblGrid>
<w:tr>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>Product</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>Price</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>Stock</w:t></w:r></w:p></w:tc>
</w:tr>
<w:tr>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>Apple</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>1.99</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>100</w:t></w:r></w:p></w:tc>
</w:tr>
<w:tr>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>Banana</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>0.99</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>50</w:t></w:r></w:p></w:tc>
</w:tr>
<w:tr>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>Orange</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>2.49</w:t></w:r></w:p></w:tc>
<w:tc><w:tcPr><w:tcW w:type='dxa' w:w='2880'/></w:tcPr><w:p><w:r><w:t>75</w:t></w:r></w:p></w:tc>
</w:tr>
</w:tbl>
"@

$rng.InsertXML($xml)
